# Weekly VCI3M forecast data refresh (NDMA pilot pages / NDVI monitoring / late NOV update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: forecast week-ending dates (serial), shifted forward one week
$dates = @(45613, 45620, 45627, 45634, 45641, 45648, 45655, 45662, 45669, 45676, 45683)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $dates[$i]
}

# Rows 2-21: updated VCI3M forecast values per region (column A / row labels unchanged)
$forecastRows = @{
    2 = @(72.90000000000001, 71.59999999999999, 70.09999999999999, 68.59999999999999, 67, 65.59999999999999, 64.3, 63.1, 62, 61, 60)
    3 = @(52.7, 53.1, 53, 52.4, 51.4, 50.1, 48.5, 46.9, 45.2, 43.7, 42.4)
    4 = @(53.2, 52.5, 51.4, 50.1, 48.7, 47.3, 46, 44.7, 43.6, 42.7, 41.9)
    5 = @(74.09999999999999, 74.40000000000001, 74.40000000000001, 74.09999999999999, 73.59999999999999, 72.90000000000001, 72.2, 71.3, 70.40000000000001, 69.5, 68.5)
    6 = @(69.59999999999999, 64.90000000000001, 60.1, 55.4, 51.1, 47.1, 43.6, 40.7, 38.3, 36.5, 35.4)
    7 = @(76.5, 72.59999999999999, 68.8, 65.3, 62.2, 59.4, 56.8, 54.4, 52.1, 49.7, 47.2)
    8 = @(80.7, 80.90000000000001, 80.8, 80.5, 80.09999999999999, 79.7, 79.5, 79.40000000000001, 79.40000000000001, 79.3, 79.2)
    9 = @(78.7, 78.09999999999999, 77.40000000000001, 76.7, 76.09999999999999, 75.7, 75.5, 75.2, 74.90000000000001, 74.40000000000001, 73.59999999999999)
    10 = @(79.8, 80, 79.40000000000001, 78, 75.90000000000001, 73.2, 70.09999999999999, 66.7, 63.2, 59.6, 56.1)
    11 = @(85.2, 85.3, 85.3, 85.09999999999999, 84.90000000000001, 84.7, 84.3, 83.8, 83.2, 82.2, 81)
    12 = @(70.8, 71.8, 72.3, 72.40000000000001, 72.2, 71.8, 71.2, 70.5, 69.8, 68.90000000000001, 67.90000000000001)
    13 = @(68.90000000000001, 68.90000000000001, 68.5, 67.7, 66.5, 65, 63.3, 61.4, 59.4, 57.4, 55.4)
    14 = @(85.59999999999999, 81.09999999999999, 76.90000000000001, 73, 69.90000000000001, 67.7, 66.3, 65.59999999999999, 65.59999999999999, 66, 66.59999999999999)
    15 = @(79.40000000000001, 78.7, 78.09999999999999, 77.8, 77.8, 78, 78.5, 79.09999999999999, 79.5, 79.7, 79.5)
    16 = @(63.6, 63.1, 62.3, 61.1, 59.6, 57.9, 56.1, 54.3, 52.5, 50.9, 49.4)
    17 = @(69.7, 69.40000000000001, 68.7, 67.8, 66.59999999999999, 65.40000000000001, 64.09999999999999, 62.8, 61.6, 60.4, 59.3)
    18 = @(74.40000000000001, 74.59999999999999, 74.7, 74.7, 74.7, 74.7, 74.7, 74.7, 74.5, 74.09999999999999, 73.5)
    19 = @(73.3, 68.3, 63.3, 58.5, 54.2, 50.6, 47.5, 45.2, 43.4, 42.1, 41.3)
    20 = @(79.5, 75.59999999999999, 72.3, 69.8, 68.2, 67.59999999999999, 67.8, 68.59999999999999, 69.59999999999999, 70.59999999999999, 71.3)
    21 = @(69.3, 67.2, 64.8, 62.2, 59.7, 57.4, 55.3, 53.6, 52.1, 51, 50.1)
}

foreach ($r in $forecastRows.Keys) {
    $vals = $forecastRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}
